$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rows = @(3, 4, 6, 7)
foreach ($r in $rows) {
    $ws.Range("E$r").Value = "[-, -, 'MCT-3A-Microcontroladores', -]"
    $ws.Range("F$r").Value = "-"
}
